# Adds a new "Search Page" test-data row plus four new columns
# (User_Name, User_Role, Employee_Name, Status) to the Test_Data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Data")

# --- New data row (row 5) -------------------------------------------------
$ws.Range("A5").Value = "searchPageTestValidData"
$ws.Range("B5").Value = "Verify Search Page"
$ws.Range("C5").Value = "Y"
$ws.Range("D5").Value = "Positive"
$ws.Range("E5").Value = "Admin"
$ws.Range("F5").Value = "admin123"
$ws.Range("G5").Value = "Admin"
$ws.Range("H5").Value = "Admin"

# --- New header columns (row 1) -------------------------------------------
$ws.Range("H1").Value = "User_Role"
$ws.Range("I1").Value = "Employee_Name"
$ws.Range("J1").Value = "Status"
$ws.Range("G1").Value = "User_Name"

# --- Remaining new-row value -----------------------------------------------
$ws.Range("J5").Value = "Enabled"

# --- Resize the new columns to fit their content, like Excel does on entry
$ws.Columns.Item(7).AutoFit() | Out-Null
$ws.Columns.Item(8).AutoFit() | Out-Null
$ws.Columns.Item(9).AutoFit() | Out-Null
$ws.Columns.Item(10).AutoFit() | Out-Null

# --- Update view / selection to match where the user left off -------------
$ws.Range("J9").Select() | Out-Null
